# Insert two new simulation-data rows after the header/first row block,
# matching the commit "Updated notebook, reran simulation".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 4-5; this pushes the former rows 4-29 down to 6-31
$ws.Rows("4:5").Insert()

# Copy the bordered/centered style used by column A down into the two new header cells
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New row 4: "Holden"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.003688183182284
$ws.Range("D4").Value = 0.9856654679821821
$ws.Range("E4").Value = 1.009370487501142
$ws.Range("F4").Value = 0.9856654679821821
$ws.Range("G4").Value = 0.9972000071007486
$ws.Range("H4").Value = 1.003525550343229
$ws.Range("I4").Value = 0.9906977083011798
$ws.Range("J4").Value = 1.009370487501142
$ws.Range("K4").Value = 1.009370487501142
$ws.Range("L4").Value = 0.9969261365792498
$ws.Range("M4").Value = 1.005331824976927
$ws.Range("N4").Value = 1.009370487501142
$ws.Range("O4").Value = 1.003688183182284
$ws.Range("P4").Value = 0.9946768255822329
$ws.Range("Q4").Value = 1.000307159880767
$ws.Range("R4").Value = 0.9995747128885361
$ws.Range("S4").Value = 0.9954265959145719
$ws.Range("T4").Value = 0.9995747128885361
$ws.Range("U4").Value = 0.9989125688112145
$ws.Range("V4").Value = 1.0010041525492
$ws.Range("W4").Value = 0.9990506707458677

# New row 5: "Rizzie Spiral"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.015172326249622
$ws.Range("D5").Value = 0.9408803744662098
$ws.Range("E5").Value = 1.0387154509201
$ws.Range("F5").Value = 0.9408803744662098
$ws.Range("G5").Value = 0.9884272230425838
$ws.Range("H5").Value = 1.0145618987606
$ws.Range("I5").Value = 0.9616280267195535
$ws.Range("J5").Value = 1.0387154509201
$ws.Range("K5").Value = 1.0387154509201
$ws.Range("L5").Value = 0.9873201960557418
$ws.Range("M5").Value = 1.021985884385018
$ws.Range("N5").Value = 1.0387154509201
$ws.Range("O5").Value = 1.015172326249622
$ws.Range("P5").Value = 0.978026350357916
$ws.Range("Q5").Value = 1.001246261152682
$ws.Range("R5").Value = 0.9982560505453106
$ws.Range("S5").Value = 0.9811242989238579
$ws.Range("T5").Value = 0.9982560505453106
$ws.Range("U5").Value = 0.9955220869229184
$ws.Range("V5").Value = 1.004160759722355
$ws.Range("W5").Value = 0.9960864225749286

# Rename category "Thomas Hex" -> "Matthies Hex" (now located at row 11 after the shift)
$ws.Range("B11").Value = "Matthies Hex"
